$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The edit is a cyclic re-shuffle of the content of rows 2,3,4,5,8,9,10
# (row 2 gets what row 8 had, row 3 gets what row 2 had, ... row 10 gets
# what row 5 had). Rows 6 and 7 are untouched. We stage every source row
# in a scratch area first so no data is overwritten before it is consumed,
# then copy staged rows into their final destination, then tidy up.
# ------------------------------------------------------------------

# 1) Stage the 7 affected rows into scratch rows far below the used range
$ws.Range("A2:AY2").Copy($ws.Range("A101:AY101"))
$ws.Range("A3:AY3").Copy($ws.Range("A102:AY102"))
$ws.Range("A4:AY4").Copy($ws.Range("A103:AY103"))
$ws.Range("A5:AY5").Copy($ws.Range("A104:AY104"))
$ws.Range("A8:AY8").Copy($ws.Range("A105:AY105"))
$ws.Range("A9:AY9").Copy($ws.Range("A106:AY106"))
$ws.Range("A10:AY10").Copy($ws.Range("A107:AY107"))

# 2) Clear the original rows
$ws.Range("A2:AY2").ClearContents()
$ws.Range("A3:AY3").ClearContents()
$ws.Range("A4:AY4").ClearContents()
$ws.Range("A5:AY5").ClearContents()
$ws.Range("A8:AY8").ClearContents()
$ws.Range("A9:AY9").ClearContents()
$ws.Range("A10:AY10").ClearContents()

# 3) Copy staged content into its new destination row
$ws.Range("A105:AY105").Copy($ws.Range("A2:AY2"))
$ws.Range("A101:AY101").Copy($ws.Range("A3:AY3"))
$ws.Range("A102:AY102").Copy($ws.Range("A4:AY4"))
$ws.Range("A103:AY103").Copy($ws.Range("A5:AY5"))
$ws.Range("A106:AY106").Copy($ws.Range("A8:AY8"))
$ws.Range("A107:AY107").Copy($ws.Range("A9:AY9"))
$ws.Range("A104:AY104").Copy($ws.Range("A10:AY10"))

# 4) Remove scratch rows
$ws.Range("A101:AY107").ClearContents()

# 5) Copying a full A:AY range also pastes blank placeholder cells into
#    columns that are structurally absent from that row; strip those back
#    out so each row keeps exactly the same cell layout as the source row.
$ws.Range("J2,K2,L2,M2,N2,O2,X2,AC2,AF2,AH2,AI2,AJ2,AK2,AL2,AM2,AN2,AO2,AP2,AQ2,AR2,AS2,AU2,AV2").ClearContents()
$ws.Range("J3,O3,X3,AF3,AH3,AI3,AJ3,AK3,AL3,AM3,AN3,AO3,AP3,AQ3,AR3,AS3,AU3,AV3").ClearContents()
$ws.Range("J4,K4,L4,M4,N4,O4,X4,AC4,AF4,AH4,AI4,AJ4,AK4,AL4,AM4,AN4,AO4,AP4,AQ4,AR4,AS4,AU4,AV4").ClearContents()
$ws.Range("J5,K5,L5,N5,O5,X5,AC5,AF5,AH5,AI5,AJ5,AK5,AL5,AM5,AN5,AO5,AP5,AQ5,AR5,AS5,AU5,AV5").ClearContents()
$ws.Range("J8,K8,L8,M8,N8,O8,X8,AC8,AF8,AH8,AI8,AJ8,AK8,AL8,AM8,AN8,AO8,AP8,AQ8,AR8,AS8,AU8,AV8").ClearContents()
$ws.Range("J9,K9,L9,M9,N9,O9,X9,AC9,AF9,AH9,AI9,AJ9,AK9,AL9,AM9,AN9,AO9,AP9,AQ9,AR9,AS9,AU9,AV9").ClearContents()
$ws.Range("J10,K10,L10,M10,N10,O10,X10,AC10,AF10,AH10,AI10,AJ10,AK10,AL10,AM10,AN10,AO10,AP10,AQ10,AR10,AS10,AU10,AV10").ClearContents()
